# Export_Tables.xlsx — "Exported regression shell table and added screenshots"
#
# Fills in the regression results (b, standard error, p-value) for the
# "Smoke" and "Constant" rows of the shell table on the
# "Regression Shell Table 1" sheet, and stamps an export timestamp a few
# rows below the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Regression Shell Table 1")

# Row 3 ("Smoke"): b, standard error, p-value
$ws.Range("B3").Value = -282.6592244418332
$ws.Range("C3").Value = 106.95441538661449
$ws.Range("D3").Value = 0.0089195273423114934

# Row 4 ("Constant"): b, standard error, p-value
$ws.Range("B4").Value = 3054.9565217391305
$ws.Range("C4").Value = 66.924275134525814
$ws.Range("D4").Value = [double]"2.5090064061084927e-103"

# Export timestamp, written further down the sheet
$ws.Range("A30").Value = "12:37:34  21 Oct 2018"
